$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 115.46154
$ws.Range("I9").Value = 91.3
$ws.Range("J9").Value = 196
$ws.Range("K9").Value = 91.3
$ws.Range("L9").Value = 196
$ws.Range("M9").Value = 77.7
$ws.Range("N9").Value = -534
$ws.Range("H40").Value = 1040
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1040
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1040
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -1390
$ws.Range("H116").Value = 2306.9167
$ws.Range("J116").Value = 2218.6365
$ws.Range("L116").Value = 2218.6365
$ws.Range("N116").Value = -9102.636500000001
$ws.Range("H133").Value = 37768.57
$ws.Range("J133").Value = 37768.57
$ws.Range("L133").Value = 37768.57
$ws.Range("N133").Value = -47888.57
$ws.Range("H137").Value = 3293.8
$ws.Range("I137").Value = 3268.7778
$ws.Range("J137").Value = 3331.3333
$ws.Range("K137").Value = 9806.3334
$ws.Range("L137").Value = 9993.999899999999
$ws.Range("M137").Value = -7256.3334
$ws.Range("N137").Value = -15093.9999
$ws.Range("H138").Value = 2735.3635
$ws.Range("J138").Value = 2686.043
$ws.Range("L138").Value = 8058.129000000001
$ws.Range("N138").Value = -18338.129

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3025.2666
$ws.Range("I74").Value = 1546
$ws.Range("K74").Value = 1546
$ws.Range("M74").Value = -672
$ws.Range("H77").Value = 3025.2666
$ws.Range("I77").Value = 1546
$ws.Range("K77").Value = 7730
$ws.Range("M77").Value = -3362
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H131").Value = 34750
$ws.Range("J131").Value = 34750
$ws.Range("L131").Value = 34750
$ws.Range("N131").Value = -44830

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 8223.429
$ws.Range("I82").Value = 4427.3335
$ws.Range("K82").Value = 4427.3335
$ws.Range("M82").Value = -4044.3335
$ws.Range("H85").Value = 8223.429
$ws.Range("I85").Value = 4427.3335
$ws.Range("K85").Value = 4427.3335
$ws.Range("M85").Value = -3101.3335
$ws.Range("H105").Value = 166667710
$ws.Range("I105").Value = 200000860
$ws.Range("K105").Value = 200000860
$ws.Range("M105").Value = -199999113
$ws.Range("H134").Value = 7758.375
$ws.Range("J134").Value = 12819.333
$ws.Range("L134").Value = 38457.999
$ws.Range("N134").Value = -43527.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1521.4807
$ws.Range("I31").Value = 1492.902
$ws.Range("J31").Value = 2979
$ws.Range("K31").Value = 1492.902
$ws.Range("L31").Value = 2979
$ws.Range("M31").Value = -1197.902
$ws.Range("N31").Value = -3569
$ws.Range("H34").Value = 1521.4807
$ws.Range("I34").Value = 1492.902
$ws.Range("J34").Value = 2979
$ws.Range("K34").Value = 1492.902
$ws.Range("L34").Value = 2979
$ws.Range("M34").Value = -1290.902
$ws.Range("N34").Value = -3383
$ws.Range("H132").Value = 1776.5946
$ws.Range("J132").Value = 2423
$ws.Range("L132").Value = 7269
$ws.Range("N132").Value = -12329
$ws.Range("H134").Value = 13515106
$ws.Range("I134").Value = 1618.44
$ws.Range("J134").Value = 41668204
$ws.Range("K134").Value = 4855.32
$ws.Range("L134").Value = 125004612
$ws.Range("M134").Value = -2320.32
$ws.Range("N134").Value = -125009682

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 2745.8
$ws.Range("J106").Value = 2745.8
$ws.Range("L106").Value = 8237.400000000001
$ws.Range("N106").Value = -10129.4
$ws.Range("H110").Value = 10562.25
$ws.Range("I110").Value = 998
$ws.Range("J110").Value = 11928.571
$ws.Range("K110").Value = 2994
$ws.Range("L110").Value = 35785.713
$ws.Range("M110").Value = 1096
$ws.Range("N110").Value = -43965.713
$ws.Range("H111").Value = 4263.1875
$ws.Range("I111").Value = 1368.5
$ws.Range("K111").Value = 4105.5
$ws.Range("M111").Value = -1038.5
$ws.Range("H112").Value = 50011040
$ws.Range("I112").Value = 4000
$ws.Range("J112").Value = 55567376
$ws.Range("K112").Value = 12000
$ws.Range("L112").Value = 166702128
$ws.Range("M112").Value = -10892
$ws.Range("N112").Value = -166704344
$ws.Range("H131").Value = 26356196
$ws.Range("I131").Value = 83333750
$ws.Range("J131").Value = 58863.77
$ws.Range("K131").Value = 250001250
$ws.Range("L131").Value = 176591.31
$ws.Range("M131").Value = -249996210
$ws.Range("N131").Value = -186671.31

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = ""
$ws.Range("H70").Value = 50003668
$ws.Range("I70").Value = 62503376
$ws.Range("J70").Value = 40003900
$ws.Range("K70").Value = 62503376
$ws.Range("L70").Value = 40003900
$ws.Range("M70").Value = -62503106
$ws.Range("N70").Value = -40004440
$ws.Range("H73").Value = 50003668
$ws.Range("I73").Value = 62503376
$ws.Range("J73").Value = 40003900
$ws.Range("K73").Value = 62503376
$ws.Range("L73").Value = 40003900
$ws.Range("M73").Value = -62502440
$ws.Range("N73").Value = -40005772
$ws.Range("H80").Value = 4185.5
$ws.Range("I80").Value = 2005
$ws.Range("K80").Value = 2005
$ws.Range("M80").Value = -1007
$ws.Range("H83").Value = 4185.5
$ws.Range("I83").Value = 2005
$ws.Range("K83").Value = 10025
$ws.Range("M83").Value = -5033
$ws.Range("H102").Value = 2863.1155
$ws.Range("I102").Value = 2068.6875
$ws.Range("J102").Value = 4134.2
$ws.Range("K102").Value = 2068.6875
$ws.Range("L102").Value = 4134.2
$ws.Range("M102").Value = -446.6875
$ws.Range("N102").Value = -7378.2
$ws.Range("H132").Value = 5580.1514
$ws.Range("I132").Value = 6909.6
$ws.Range("J132").Value = 3534.8462
$ws.Range("K132").Value = 20728.8
$ws.Range("L132").Value = 10604.5386
$ws.Range("M132").Value = -18198.8
$ws.Range("N132").Value = -15664.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2497.5
$ws.Range("I7").Value = 2490
$ws.Range("J7").Value = 2501.25
$ws.Range("K7").Value = 2490
$ws.Range("L7").Value = 2501.25
$ws.Range("M7").Value = -2378
$ws.Range("N7").Value = -2725.25
$ws.Range("H13").Value = 2985.7144
$ws.Range("I13").Value = 2985.7144
$ws.Range("K13").Value = 2985.7144
$ws.Range("M13").Value = -2845.7144
$ws.Range("H22").Value = 1333.3334
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 1333.3334
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1714
$ws.Range("H46").Value = 4153.636
$ws.Range("I46").Value = 100
$ws.Range("J46").Value = 4559
$ws.Range("K46").Value = 100
$ws.Range("L46").Value = 4559
$ws.Range("M46").Value = 88
$ws.Range("N46").Value = -4935
$ws.Range("H122").Value = 31252938
$ws.Range("I122").Value = 41669416
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 125008248
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -125005798
$ws.Range("N122").Value = -15400
$ws.Range("H126").Value = 2497.5
$ws.Range("I126").Value = 2490
$ws.Range("J126").Value = 2501.25
$ws.Range("K126").Value = 7470
$ws.Range("L126").Value = 7503.75
$ws.Range("M126").Value = -5000
$ws.Range("N126").Value = -12443.75
$ws.Range("H132").Value = 2504.6553
$ws.Range("I132").Value = 2125.1333
$ws.Range("J132").Value = 2911.2856
$ws.Range("K132").Value = 6375.3999
$ws.Range("L132").Value = 8733.856800000001
$ws.Range("M132").Value = -3845.3999
$ws.Range("N132").Value = -13793.8568
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""
$ws.Range("H136").Value = 1831.8422
$ws.Range("I136").Value = 1766.9445
$ws.Range("K136").Value = 5300.833500000001
$ws.Range("M136").Value = -2750.833500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1175
$ws.Range("I96").Value = 1824.75
$ws.Range("J96").Value = 741.8333
$ws.Range("K96").Value = 1824.75
$ws.Range("L96").Value = 741.8333
$ws.Range("M96").Value = -451.75
$ws.Range("N96").Value = -3487.8333
$ws.Range("H100").Value = 288.84616
$ws.Range("I100").Value = 312.44446
$ws.Range("J100").Value = 235.75
$ws.Range("K100").Value = 624.88892
$ws.Range("L100").Value = 471.5
$ws.Range("M100").Value = -83.88891999999998
$ws.Range("N100").Value = -1553.5
$ws.Range("H126").Value = 90910320
$ws.Range("I126").Value = 142858420
$ws.Range("K126").Value = 428575260
$ws.Range("M126").Value = -428572790
$ws.Range("H132").Value = 3292.4583
$ws.Range("I132").Value = 3441.923
$ws.Range("J132").Value = 2644.7778
$ws.Range("K132").Value = 10325.769
$ws.Range("L132").Value = 7934.3334
$ws.Range("M132").Value = -7795.769
$ws.Range("N132").Value = -12994.3334
